$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (formulas text) into column B
$ws.Range("A1").EntireColumn.Insert()

# Set the new column A values (cluster labels)
$ws.Range("A1").Value = "IgGI"
$ws.Range("A2").Value = "IgGII"

# Set column B width to match target (~14.75 characters)
$ws.Range("B1").EntireColumn.ColumnWidth = 14

# Update selection to match target state
$ws.Range("C8").Select()
